$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 6) - B6 entered before A6 to match shared string order
$ws.Range("B6").Value = "validate after appointment return to homepage"
$ws.Range("A6").Value = "Tc_005"

# Set row height for the new row
$ws.Rows.Item(6).RowHeight = 15.6

# Apply font/alignment formatting to B6: size 12, vertical centered
$ws.Range("B6").Font.Size = 12
$ws.Range("B6").VerticalAlignment = -4108

# Widen column B to fit the new content (target stored width 49.5546875 chars)
$ws.Columns.Item(2).ColumnWidth = 48.7213541666667

# Update the selected cell to match the target workbook view
$ws.Range("F13").Select()
